$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 - text columns (A-E)
$ws.Range("A18").Value = "Brazilian Serie A"
# B18 holds a literal "yyyy-mm-dd" string (like the other Date cells in this
# sheet) rather than a real date, so force text formatting before assigning
# it to stop the auto date-detection from turning it into a date serial.
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2025-12-02"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "21:30:00"
$ws.Range("D18").Value = "Gremio"
$ws.Range("E18").Value = "Fluminense"

# Updated / new numeric cell values
$ws.Range("I3").Value = 2.38
$ws.Range("L3").Value = 1.37
$ws.Range("M3").Value = 1.09
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 1.39
$ws.Range("R3").Value = 1.26
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 1.88
$ws.Range("U3").Value = 1.93
$ws.Range("X3").Value = 14
$ws.Range("Y3").Value = 10.5
$ws.Range("Z3").Value = 16.5
$ws.Range("AA3").Value = 38
$ws.Range("AB3").Value = 15
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AD3").Value = 13
$ws.Range("AE3").Value = 34
$ws.Range("AF3").Value = 34
$ws.Range("AG3").Value = 19.5
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 60
$ws.Range("AJ3").Value = 100
$ws.Range("AK3").Value = 65
$ws.Range("AL3").Value = 80
$ws.Range("AM3").Value = 160
$ws.Range("AN3").Value = 80
$ws.Range("AO3").Value = 28
$ws.Range("G4").Value = 2.2
$ws.Range("L4").Value = 1.43
$ws.Range("W4").Value = 1.83
$ws.Range("Z4").Value = 29
$ws.Range("F5").Value = 2.06
$ws.Range("I5").Value = 4.3
$ws.Range("J5").Value = 3.2
$ws.Range("U5").Value = 1.94
$ws.Range("V5").Value = 1.31
$ws.Range("AO5").Value = 80
$ws.Range("P7").Value = 1.71
$ws.Range("Q7").Value = 1.96
$ws.Range("H8").Value = 3.3
$ws.Range("F9").Value = 5.5
$ws.Range("L9").Value = 1.23
$ws.Range("N9").Value = 2.32
$ws.Range("P9").Value = 2.06
$ws.Range("R9").Value = 1.44
$ws.Range("S9").Value = 2.78
$ws.Range("T9").Value = 1.81
$ws.Range("U9").Value = 1.98
$ws.Range("X9").Value = 20
$ws.Range("Y9").Value = 11
$ws.Range("Z9").Value = 12
$ws.Range("AA9").Value = 17.5
$ws.Range("AB9").Value = 30
$ws.Range("AC9").Value = 12.5
$ws.Range("AD9").Value = 12.5
$ws.Range("AE9").Value = 19
$ws.Range("AF9").Value = 70
$ws.Range("AG9").Value = 32
$ws.Range("AH9").Value = 26
$ws.Range("AI9").Value = 40
$ws.Range("AJ9").Value = 230
$ws.Range("AK9").Value = 120
$ws.Range("AL9").Value = 110
$ws.Range("AM9").Value = 140
$ws.Range("AN9").Value = 130
$ws.Range("AO9").Value = 9.199999999999999
$ws.Range("F10").Value = 2.18
$ws.Range("L11").Value = 1.32
$ws.Range("N11").Value = 5
$ws.Range("M12").Value = 1.08
$ws.Range("X12").Value = 13
$ws.Range("F14").Value = 1.75
$ws.Range("AG14").Value = 9.800000000000001
$ws.Range("F15").Value = 1.9
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 4.7
$ws.Range("I15").Value = 5.3
$ws.Range("J15").Value = 3.3
$ws.Range("K15").Value = 3.7
$ws.Range("L15").Value = 1.48
$ws.Range("N15").Value = 3.3
$ws.Range("P15").Value = 1.78
$ws.Range("Q15").Value = 2.16
$ws.Range("R15").Value = 1.28
$ws.Range("W15").Value = 2
$ws.Range("X15").Value = 14
$ws.Range("F17").Value = 2.08
$ws.Range("J17").Value = 3.65
$ws.Range("K17").Value = 3.95
$ws.Range("Q17").Value = 1.87
$ws.Range("R17").Value = 1.37
$ws.Range("W17").Value = 1.81
$ws.Range("Y17").Value = 980
$ws.Range("F18").Value = 2.96
$ws.Range("G18").Value = 3.3
$ws.Range("H18").Value = 2.56
$ws.Range("I18").Value = 2.8
$ws.Range("J18").Value = 3.15
$ws.Range("K18").Value = 3.45
$ws.Range("L18").Value = 1.49
$ws.Range("M18").Value = 1.1
$ws.Range("N18").Value = 2.84
$ws.Range("O18").Value = 1.41
$ws.Range("P18").Value = 1.67
$ws.Range("Q18").Value = 2.28
$ws.Range("R18").Value = 1.24
$ws.Range("S18").Value = 4.2
$ws.Range("T18").Value = 1.89
$ws.Range("U18").Value = 1.93
$ws.Range("V18").Value = 1.56
$ws.Range("W18").Value = 1.44
$ws.Range("X18").Value = 11.5
$ws.Range("Y18").Value = 9.6
$ws.Range("Z18").Value = 17.5
$ws.Range("AA18").Value = 980
$ws.Range("AB18").Value = 11
$ws.Range("AC18").Value = 7.8
$ws.Range("AD18").Value = 13
$ws.Range("AE18").Value = 980
$ws.Range("AF18").Value = 22
$ws.Range("AG18").Value = 14.5
$ws.Range("AH18").Value = 21
$ws.Range("AI18").Value = 55
$ws.Range("AJ18").Value = 60
$ws.Range("AK18").Value = 980
$ws.Range("AL18").Value = 60
$ws.Range("AM18").Value = 160
$ws.Range("AN18").Value = 980
$ws.Range("AO18").Value = 980
